$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..42 down to 4..43.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new data record.
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(3, 3).Value = 'Maule'
$ws.Cells.Item(3, 4).Value = 44685
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = 'Fruta'
$ws.Cells.Item(3, 7).Value = 100104
$ws.Cells.Item(3, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(3, 9).Value = 100104003
$ws.Cells.Item(3, 10).Value = 'Membrillo'
$ws.Cells.Item(3, 11).Value = 'Champion'
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 190
$ws.Cells.Item(3, 14).Value = 10000
$ws.Cells.Item(3, 15).Value = 10000
$ws.Cells.Item(3, 16).Value = 10000
$ws.Cells.Item(3, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(3, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value = 556
$ws.Cells.Item(3, 20).Value = 18
